$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G7").Value = "my new change"
$ws.Range("G7").Select()
